# Scheduled-runner style refresh of market-price/profit columns (H:N) on
# several Leve sheets (currentAveragePrice*, LevePrice*, LeveProfit*).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 302.4
$ws.Range("I12").Value = 359.25
$ws.Range("J12").Value = 75
$ws.Range("K12").Value = 359.25
$ws.Range("L12").Value = 75
$ws.Range("M12").Value = -189.25
$ws.Range("N12").Value = -415
$ws.Range("H92").Value = 3007.8572
$ws.Range("I92").Value = 3143.25
$ws.Range("J92").Value = 300
$ws.Range("K92").Value = 3143.25
$ws.Range("L92").Value = 300
$ws.Range("M92").Value = -1895.25
$ws.Range("N92").Value = -2796
$ws.Range("H95").Value = 29832.334
$ws.Range("J95").Value = 29832.334
$ws.Range("L95").Value = 29832.334
$ws.Range("N95").Value = -35324.334
$ws.Range("H97").Value = 2750
$ws.Range("J97").Value = 2750
$ws.Range("L97").Value = 8250
$ws.Range("N97").Value = -9242
$ws.Range("H118").Value = 91654.37
$ws.Range("I118").Value = 143527
$ws.Range("J118").Value = 877.25
$ws.Range("K118").Value = 430581
$ws.Range("L118").Value = 2631.75
$ws.Range("M118").Value = -428924
$ws.Range("N118").Value = -5945.75
$ws.Range("H132").Value = 3134.868
$ws.Range("I132").Value = 1974.4722
$ws.Range("J132").Value = 5592.1763
$ws.Range("K132").Value = 5923.4166
$ws.Range("L132").Value = 16776.5289
$ws.Range("M132").Value = -3393.4166
$ws.Range("N132").Value = -21836.5289
$ws.Range("H135").Value = 527
$ws.Range("I135").Value = 341.75
$ws.Range("J135").Value = 2750
$ws.Range("K135").Value = 3075.75
$ws.Range("L135").Value = 24750
$ws.Range("M135").Value = -540.75
$ws.Range("N135").Value = -29820
$ws.Range("H137").Value = 3541.658
$ws.Range("I137").Value = 3419.2693
$ws.Range("K137").Value = 10257.8079
$ws.Range("M137").Value = -7707.8079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6654.791
$ws.Range("I32").Value = 4284.8535
$ws.Range("J32").Value = 28247.555
$ws.Range("K32").Value = 4284.8535
$ws.Range("L32").Value = 28247.555
$ws.Range("M32").Value = -3997.8535
$ws.Range("N32").Value = -28821.555
$ws.Range("H61").Value = 1780.39
$ws.Range("I61").Value = 1064.6349
$ws.Range("J61").Value = 2999.1082
$ws.Range("K61").Value = 1064.6349
$ws.Range("L61").Value = 2999.1082
$ws.Range("M61").Value = -852.6349
$ws.Range("N61").Value = -3423.1082
$ws.Range("H132").Value = 2678.1667
$ws.Range("J132").Value = 2433.4614
$ws.Range("L132").Value = 7300.3842
$ws.Range("N132").Value = -12360.3842
$ws.Range("H136").Value = 1780.39
$ws.Range("I136").Value = 1064.6349
$ws.Range("J136").Value = 2999.1082
$ws.Range("K136").Value = 3193.9047
$ws.Range("L136").Value = 8997.3246
$ws.Range("M136").Value = -643.9047
$ws.Range("N136").Value = -14097.3246

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 754.6061
$ws.Range("I64").Value = 699.4545000000001
$ws.Range("K64").Value = 699.4545000000001
$ws.Range("M64").Value = -474.4545000000001
$ws.Range("H67").Value = 754.6061
$ws.Range("I67").Value = 699.4545000000001
$ws.Range("K67").Value = 699.4545000000001
$ws.Range("M67").Value = 80.54549999999995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2291.9866
$ws.Range("I31").Value = 1372.3405
$ws.Range("J31").Value = 3892.8518
$ws.Range("K31").Value = 1372.3405
$ws.Range("L31").Value = 3892.8518
$ws.Range("M31").Value = -1077.3405
$ws.Range("N31").Value = -4482.8518
$ws.Range("H34").Value = 2291.9866
$ws.Range("I34").Value = 1372.3405
$ws.Range("J34").Value = 3892.8518
$ws.Range("K34").Value = 1372.3405
$ws.Range("L34").Value = 3892.8518
$ws.Range("M34").Value = -1170.3405
$ws.Range("N34").Value = -4296.8518
$ws.Range("H58").Value = 2718.4062
$ws.Range("I58").Value = 2957.318
$ws.Range("J58").Value = 2192.8
$ws.Range("K58").Value = 2957.318
$ws.Range("L58").Value = 2192.8
$ws.Range("M58").Value = -2754.318
$ws.Range("N58").Value = -2598.8
$ws.Range("H132").Value = 1647.9565
$ws.Range("I132").Value = 874.42426
$ws.Range("J132").Value = 3611.5386
$ws.Range("K132").Value = 2623.27278
$ws.Range("L132").Value = 10834.6158
$ws.Range("M132").Value = -93.27278000000024
$ws.Range("N132").Value = -15894.6158
$ws.Range("H134").Value = 1484.52
$ws.Range("I134").Value = 934.68085
$ws.Range("J134").Value = 2407.4644
$ws.Range("K134").Value = 2804.04255
$ws.Range("L134").Value = 7222.3932
$ws.Range("M134").Value = -269.0425500000001
$ws.Range("N134").Value = -12292.3932
$ws.Range("H136").Value = 2718.4062
$ws.Range("I136").Value = 2957.318
$ws.Range("J136").Value = 2192.8
$ws.Range("K136").Value = 8871.954000000002
$ws.Range("L136").Value = 6578.400000000001
$ws.Range("M136").Value = -6321.954000000002
$ws.Range("N136").Value = -11678.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 582.6316
$ws.Range("I5").Value = 476.48148
$ws.Range("J5").Value = 843.1818
$ws.Range("K5").Value = 1429.44444
$ws.Range("L5").Value = 2529.5454
$ws.Range("M5").Value = -1317.44444
$ws.Range("N5").Value = -2753.5454
$ws.Range("H92").Value = 567.9
$ws.Range("J92").Value = 554.1429000000001
$ws.Range("L92").Value = 1662.4287
$ws.Range("N92").Value = -4158.4287
$ws.Range("H93").Value = 1675
$ws.Range("J93").Value = 1566.6666
$ws.Range("L93").Value = 4699.9998
$ws.Range("N93").Value = -8443.9998
$ws.Range("H122").Value = 1024.9
$ws.Range("I122").Value = 574.75
$ws.Range("J122").Value = 1325
$ws.Range("K122").Value = 5172.75
$ws.Range("L122").Value = 11925
$ws.Range("M122").Value = -2722.75
$ws.Range("N122").Value = -16825
$ws.Range("H135").Value = 582.6316
$ws.Range("I135").Value = 476.48148
$ws.Range("J135").Value = 843.1818
$ws.Range("K135").Value = 4288.33332
$ws.Range("L135").Value = 7588.6362
$ws.Range("M135").Value = -1753.33332
$ws.Range("N135").Value = -12658.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3810.825
$ws.Range("I132").Value = 3624
$ws.Range("J132").Value = 4091.0625
$ws.Range("K132").Value = 10872
$ws.Range("L132").Value = 12273.1875
$ws.Range("M132").Value = -8342
$ws.Range("N132").Value = -17333.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5429.85
$ws.Range("I132").Value = 2274
$ws.Range("J132").Value = 10506.652
$ws.Range("K132").Value = 6822
$ws.Range("L132").Value = 31519.956
$ws.Range("M132").Value = -4292
$ws.Range("N132").Value = -36579.956

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1711.6323
$ws.Range("I132").Value = 1242.9535
$ws.Range("J132").Value = 2517.76
$ws.Range("K132").Value = 3728.8605
$ws.Range("L132").Value = 7553.280000000001
$ws.Range("M132").Value = -1198.8605
$ws.Range("N132").Value = -12613.28
$ws.Range("H136").Value = 11757454
$ws.Range("I136").Value = 18888188
$ws.Range("J136").Value = 305064.25
$ws.Range("K136").Value = 56664564
$ws.Range("L136").Value = 915192.75
$ws.Range("M136").Value = -56662014
$ws.Range("N136").Value = -920292.75
